$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.344.64'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.076.04'
$ws.Range("E3").Value = '  +4.52%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Formula = "'234.61"
$ws.Range("E5").Value = '  -2.91%  '
$ws.Range("D6").Formula = "'0.616"
$ws.Range("E6").Value = '  +2.21%  '
$ws.Range("D8").Formula = "'57.36"
$ws.Range("E8").Value = '  +4.19%  '
$ws.Range("D9").Formula = "'0.381"
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").Formula = "'0.101"
$ws.Range("E12").Value = '  +3.79%  '
$ws.Range("D13").Value = '2.382.54'
$ws.Range("E13").Value = '  +4.73%  '
$ws.Range("D14").Formula = "'14.49"
$ws.Range("E14").Value = '  +3.19%  '
$ws.Range("D15").Formula = "'21.01"
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Formula = "'0.774"
$ws.Range("E16").Value = '  +2.32%  '
$ws.Range("D17").Formula = "'5.23"
$ws.Range("E17").Value = '  +3.87%  '
$ws.Range("D18").Value = '2.030.35'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("D19").Value = '37.481.37'
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").Formula = "'6.00"
$ws.Range("E20").Value = '  +20.72%  '
$ws.Range("D21").Formula = "'68.31"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Formula = "'223.29"
$ws.Range("E23").Value = '  -1.95%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  +3.18%  '
$ws.Range("D26").Formula = "'2.42"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").Formula = "'162.69"
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").Formula = "'8.87"
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("D29").Formula = "'0.131"
$ws.Range("D30").Formula = "'19.24"
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("D33").Formula = "'4.45"
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").Formula = "'0.0622"
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  +9.67%  '
$ws.Range("E36").Value = '  +4.56%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Formula = "'5.97"
$ws.Range("E38").Value = '  +14.26%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Formula = "'1.77"
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Formula = "'3.31"
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("E41").Value = '  -4.45%  '
$ws.Range("D42").Formula = "'0.0959"
$ws.Range("E42").Value = '  +9.34%  '
$ws.Range("D43").Value = '1.471.76'
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("D44").Formula = "'4.36"
$ws.Range("E44").Value = '  +19.01%  '
$ws.Range("D45").Formula = "'95.05"
$ws.Range("E45").Value = '  +7.87%  '
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("D47").Formula = "'16.14"
$ws.Range("E47").Value = '  +6.34%  '
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("E49").Value = '  +3.14%  '
$ws.Range("D50").Formula = "'7.28"
$ws.Range("E50").Value = '  +9.57%  '
$ws.Range("E51").Value = '  +1.76%  '
